$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "21_01_2024"
$ws.Range("E2").Value = 916
$ws.Range("E3").Value = 890
$ws.Range("E4").Value = 1068
$ws.Range("E5").Value = 2845

$ws.Range("E3").Select()
